# Insert a new data row at row 175 (pushes existing rows 175-195 down to 176-196)
# and populate it with the new "Cebollín" price record.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(175).Insert()

$ws.Range("A175").Value = 7
$ws.Range("B175").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C175").Value = "Ñuble"
$ws.Range("D175").Value = 45154
$ws.Range("E175").Value = 16
$ws.Range("F175").Value = 100112037
$ws.Range("G175").Value = "Cebollín"
$ws.Range("H175").Value = "Sin especificar"
$ws.Range("I175").Value = "Primera"
$ws.Range("J175").Value = 150
$ws.Range("K175").Value = 6000
$ws.Range("L175").Value = 6000
$ws.Range("M175").Value = 6000
$ws.Range("N175").Value = "`$/paquete 36 unidades"
$ws.Range("O175").Value = "Provincia de Diguillín"
$ws.Range("P175").Value = 167
$ws.Range("Q175").Value = 36
$ws.Range("R175").Value = "Hortaliza"
